# Commit: "Update BGDP in other files that rely on those data"
#
# About!A33 holds the India:US GDP-per-capita adjustment factor ("BGDP").
# Bumping it ripples (via formulas like `=Data!C3*About!$A$33`) into every
# dependent cell on the DCpUC sheet, so a single write here reproduces the
# whole cascade of recalculated values.

$wb = $excel.ActiveWorkbook

$aboutSheet = $wb.Worksheets.Item("About")
$aboutSheet.Range("A33").Value = 0.03878298458735905

# Mirror the reviewer's on-screen state after editing A33: the About sheet
# stays the active/selected sheet, scrolled down so row 22 is at the top,
# with A33 (the cell just edited) now the active selection.
$aboutSheet.Activate()
$excel.ActiveWindow.ScrollRow = 22
$aboutSheet.Range("A33").Select()
